$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Force Excel to store the value as literal text (matching the
    # original inlineStr cell) instead of silently parsing numeric-
    # looking strings ("27.139.25", "312.50", ...) into numbers.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Restore the cell style to the workbook default (no explicit
    # number format / quote-prefix flag), matching the original file.
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '27.139.25'
Set-TextValue 'E2' '  -0.24%  '

# Row 3
Set-TextValue 'D3' '1.826.37'
Set-TextValue 'E3' '  -0.47%  '

# Row 4
Set-TextValue 'E4' '  -0.08%  '

# Row 5
Set-TextValue 'D5' '312.50'
Set-TextValue 'E5' '  -0.39%  '

# Row 6
Set-TextValue 'E6' '  -0.18%  '

# Row 7
Set-TextValue 'D7' '0.4629'
Set-TextValue 'E7' '  -1.84%  '

# Row 8
Set-TextValue 'D8' '0.3631'
Set-TextValue 'E8' '  -1.57%  '

# Row 9
Set-TextValue 'D9' '0.07286'
Set-TextValue 'E9' '  -1.90%  '

# Row 10
Set-TextValue 'D10' '0.8692'
Set-TextValue 'E10' '  -1.56%  '

# Row 11
Set-TextValue 'D11' '20.12'
Set-TextValue 'E11' '  -1.56%  '

# Row 12
Set-TextValue 'D12' '1.879.46'
Set-TextValue 'E12' '  +3.18%  '

# Row 13
Set-TextValue 'D13' '0.07637'
Set-TextValue 'E13' '  +4.07%  '

# Row 14
Set-TextValue 'E14' '  -2.57%  '

# Row 15
Set-TextValue 'D15' '92.37'
Set-TextValue 'E15' '  -0.61%  '

# Row 16
Set-TextValue 'D16' '6.469'
Set-TextValue 'E16' '  -1.53%  '

# Row 17
Set-TextValue 'E17' '  -0.19%  '

# Row 18
Set-TextValue 'D18' '0.000008614'
Set-TextValue 'E18' '  -2.04%  '

# Row 19
Set-TextValue 'E19' '  -0.04%  '

# Row 20
Set-TextValue 'D20' '27.443.48'
Set-TextValue 'E20' '  +0.85%  '

# Row 21
Set-TextValue 'D21' '14.48'
Set-TextValue 'E21' '  -2.07%  '

# Row 22
Set-TextValue 'D22' '5.214'
Set-TextValue 'E22' '  -1.77%  '

# Row 23
Set-TextValue 'E23' '  -1.18%  '

# Row 24
Set-TextValue 'D24' '2.096.89'
Set-TextValue 'E24' '  +2.14%  '

# Row 25
Set-TextValue 'D25' '1.887'
Set-TextValue 'E25' '  -0.92%  '

# Row 26
Set-TextValue 'D26' '151.06'
Set-TextValue 'E26' '  -0.98%  '

# Row 27
Set-TextValue 'D27' '18.28'
Set-TextValue 'E27' '  -1.91%  '

# Row 28
Set-TextValue 'D28' '2.094'
Set-TextValue 'E28' '  -3.14%  '

# Row 29
Set-TextValue 'D29' '5.122'
Set-TextValue 'E29' '  -3.03%  '

# Row 30
Set-TextValue 'D30' '116.25'
Set-TextValue 'E30' '  -1.37%  '

# Row 31
Set-TextValue 'D31' '0.08911'
Set-TextValue 'E31' '  -0.25%  '

# Row 32
Set-TextValue 'D32' '2.958'
Set-TextValue 'E32' '  +0.52%  '

# Row 33
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.7365'
Set-TextValue 'E33' '  -3.17%  '

# Row 34
Set-TextValue 'B34' 'ARBITRUM'
Set-TextValue 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.151'
Set-TextValue 'E34' '  -1.88%  '

# Row 35
Set-TextValue 'D35' '4.457'
Set-TextValue 'E35' '  -2.09%  '

# Row 36
Set-TextValue 'D36' '1.011'
Set-TextValue 'E36' '  -0.08%  '

# Row 37
Set-TextValue 'D37' '2.516'
Set-TextValue 'E37' '  +3.57%  '

# Row 38
Set-TextValue 'D38' '1.082'
Set-TextValue 'E38' '  -1.91%  '

# Row 39
Set-TextValue 'D39' '0.05242'
Set-TextValue 'E39' '  -1.88%  '

# Row 40
Set-TextValue 'D40' '0.01911'
Set-TextValue 'E40' '  -2.63%  '

# Row 41
Set-TextValue 'D41' '2.926'
Set-TextValue 'E41' '  -2.70%  '

# Row 42
Set-TextValue 'D42' '7.175'
Set-TextValue 'E42' '  -2.41%  '

# Row 43
Set-TextValue 'D43' '0.5208'
Set-TextValue 'E43' '  -2.67%  '

# Row 44
Set-TextValue 'D44' '0.1629'
Set-TextValue 'E44' '  -2.10%  '

# Row 45
Set-TextValue 'D45' '8.293'
Set-TextValue 'E45' '  -2.97%  '

# Row 46
Set-TextValue 'E46' '  -2.09%  '

# Row 47
Set-TextValue 'D47' '10.27'
Set-TextValue 'E47' '  -2.34%  '

# Row 48
Set-TextValue 'D48' '1.011'

# Row 49
Set-TextValue 'D49' '103.56'
Set-TextValue 'E49' '  -0.37%  '

# Row 50
Set-TextValue 'D50' '1.637'
Set-TextValue 'E50' '  -2.13%  '

# Row 51
Set-TextValue 'D51' '0.06269'
Set-TextValue 'E51' '  -0.99%  '
